$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the previous "last sale" label from VENDA 08 (29/09) to VENDA 09 (29/09)
$ws.Range("G13").Value = "VENDA 09 (29/09)"

# Copy the formatting of the row above (row 13) down into the new row 14
# so the new client entry keeps the same borders/fill/font/number formats.
$ws.Range("B13:G13").Copy()
$ws.Range("B14:G14").PasteSpecial(-4122)

# Fill in the new client's data for row 14
$ws.Range("B14").Value = "ISAQUE GOMES"
$ws.Range("C14").Value = "e15c263fbb0290b7886838113821ff8b"
$ws.Range("D14").Value = 44833
$ws.Range("E14").Value = 365
$ws.Range("F14").Value = "-"
$ws.Range("G14").Value = "VENDA 10 (29/09)"
